# Update values in column E (imputed results) to reflect the latest
# RandomForest algorithm run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E9"  = 17.06250000000001
    "E13" = 16.6825
    "E16" = 16.669
    "E18" = 17.63940000000002
    "E20" = 16.0187
    "E26" = 16.17059999999999
    "E27" = 16.75619999999999
    "E29" = 16.86840000000002
    "E35" = 16.0135
    "E36" = 17.52920000000001
    "E45" = 16.4628
    "E55" = 16.56019999999999
    "E57" = 16.6039
    "E69" = 17.22590000000002
    "E76" = 16.21309999999999
    "E78" = 16.67800000000002
    "E82" = 16.84810000000001
    "E83" = 16.4686
    "E93" = 18.01510000000002
    "E97" = 16.592
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
